# Actualización automática hashcode mar oct  8 01:51:51 CEST 2019
# Updates specific hashcode values (column B) in the active worksheet
# to match the new values produced by the upstream hashcode generator.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    11  = "e1d9e5886ee585ad0a9b7601da0b9e2d"
    94  = "3976bbb9f1b4a382bc87fb541bc59088"
    95  = "5c73882c3c53b385a9b6cb3418168164"
    98  = "dfb77a4ff63d5cca57d0b52f6e8ac4ad"
    109 = "be3773578f4667dab4f496d33a85eeab"
    115 = "662197525b2acd21c6124243032fe1bd"
    126 = "20001017abcc4f44c78a0dd45a0cb1a0"
    159 = "e156324346350ccc26b299ced439b2e5"
    169 = "574e5e040d38a98dd0601262e0d8c5ee"
    227 = "0bed25d524905a11baaf024e5fd8abc9"
    232 = "72e804d3ceaaf08953cc162b25b3431f"
    302 = "eb82f80ba4867ba6d0317941901c73a3"
    339 = "06373edd62c8e2b9ae8898478ade4978"
    420 = "0841f66eec1f7caf51680bed6f5054c6"
    464 = "bd4dd24afc8c843a92b91727d260b2e9"
    483 = "0a8277e209a3872254017c3c03014112"
    506 = "ea1655feed1aab34539739f88df3e2a3"
    524 = "3fcdff4b97f2ae0d13fa819c55d00493"
    558 = "4d73a6bc8ae7518487d96671ff690408"
    600 = "eef16b95de2fdd043b7a987a50adf02f"
    604 = "58f975222f3a135fb1d18e6d14e8dbbe"
    624 = "8087a7ff768fba1c6fb773965ee470d9"
    626 = "124054d4a6a4cbe2c5a28c761a12800c"
    635 = "17f107c3ec809afa64d7dd72684ac46a"
    657 = "37f2fc6771073b8a151c9c9ef80f910d"
    663 = "250e948aef42581b3b154890646c7481"
    827 = "5c5abd9c802a6043120d7bf33bf913e3"
    870 = "e6fd6e055f6922f9c72c3dd36ca75b9e"
    882 = "d878f735a89572d2273c1e98708e28dd"
}

foreach ($row in $updates.Keys) {
    $ws.Cells.Item($row, 2).Value = $updates[$row]
}
